$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New rows of data for the special price kiosk (film showings 11.-13.10.24)
$rows = @(
    @(45576, "Spez 1", "Raclette klein",      6, 12),
    @(45576, "Spez 2", "Raclette gross",     12, 19),
    @(45576, "Spez 3", "Weisswein Wallis",    6, 35),
    @(45577, "Spez 1", "Raclette klein",      6, 13),
    @(45577, "Spez 2", "Raclette gross",     12,  9),
    @(45577, "Spez 3", "Weisswein Wallis",    6, 26),
    @(45578, "Spez 1", "Zauberstab-Spiessli", 4,  9)
)

$startRow = 20
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $prev = $r - 1

    # Copy formatting (number formats/styles) down from the row above
    $ws.Range("A$prev`:E$prev").Copy() | Out-Null
    $ws.Range("A$r`:E$r").PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
}
$excel.CutCopyMode = 0

# Resize the table to include the new rows
$tbl = $ws.ListObjects.Item("Table1")
$tbl.Resize($ws.Range("A1:E26"))

# Update selection to reflect last edited cell
$ws.Range("A26").Select()
